# train human model v2
# Updates the last existing row (137) with a new record and appends 16
# brand-new rows (138-153) of human/computer score data, extending the
# sheet's used range from A1:D137 to A1:D153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (A firebase-style key, computer_score, human_score, behavior)
$rows = @(
    @(137, "-N7S6teeioJOgwOkGHeb", 54, 59, "follow_stag"),
    @(138, "-N7S6zxhC8D9WPKQ1GTJ", 25, 24, "closest"),
    @(139, "-N7S7QCaT2R-cN-xQoen", 5,  8,  "random"),
    @(140, "-N7S8YwrV6R9W4BW1dzL", 9,  25, "follow_stag"),
    @(141, "-N7S92Id30UECmsUw1ZJ", 62, 61, "follow_stag"),
    @(142, "-N7SAbmiI0l11nqFTIxY", 6,  31, "follow_stag"),
    @(143, "-N7SDLpwrkF3J6tlEbbg", 42, 48, "follow_stag"),
    @(144, "-N7SE7P1i5pBPagbU2fs", 7,  22, "random"),
    @(145, "-N7SG-gcyIrxsm71UJ70", 4,  19, "follow_stag"),
    @(146, "-N7SKdAU00pOm8Cdx2lf", 5,  22, "random"),
    @(147, "-N7SNi9HqPWM2I4geMxR", 8,  24, "random"),
    @(148, "-N7SQQHk7nWiuDwke8jf", 25, 19, "closest"),
    @(149, "-N7SS0ZNRRzABEFqOvRm", 37, 51, "follow_stag"),
    @(150, "-N7SWDA_ldwYcJkG2XKL", 57, 57, "follow_stag"),
    @(151, "-N7SZPE3oHGDDzPTrYxV", 7,  19, "random"),
    @(152, "-N7ScVpejC0Dy0ua8k1T", 56, 55, "follow_stag"),
    @(153, "-N7SuS0uDb2PVST69Iae", 5,  25, "follow_stag")
)

$templateRow = 137

foreach ($r in $rows) {
    $rowNum = $r[0]

    # New rows need the same look (border box, bold, centered text) as the
    # existing data rows; cloning the template row's formatting is the
    # simplest way to keep every new row visually identical to row 137.
    if ($rowNum -ne $templateRow) {
        $ws.Range("A$templateRow`:D$templateRow").Copy($ws.Range("A$rowNum`:D$rowNum"))
    }

    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
